$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "aggregate_id"
for ($r = 2; $r -le 15; $r++) {
    $ws.Range("F$r").Value = $ws.Range("A$r").Value2
}

$ws.Range("H12").Select()
